# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders the "Periodo Mora" detail rows (16-37) on Hoja1 so the new
# account-statement data block (ELVIS CANTILLO CASTILLA) comes first,
# followed by LUZ ESTELA MARTINEZ LONDOÑO, then LEDYS DEL SOCORRO ROMERO
# PATERNINA - each worker's periods now in ascending order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data, in row order (row 16 .. row 37):
#   C = Tipo Doc Trabajador / N Doc Trabajador, D = Nombre Trabajador,
#   E = Periodo Mora, F = Valor Mora, G = Salario Basico
$rows = @(
    @("73095460", "ELVIS CANTILLO CASTILLA",            "2102", 64682,  3828128),
    @("73095460", "ELVIS CANTILLO CASTILLA",            "2103", 129365, 3828128),
    @("73095460", "ELVIS CANTILLO CASTILLA",            "2104", 129365, 3828128),
    @("73095460", "ELVIS CANTILLO CASTILLA",            "2105", 129365, 3828128),
    @("73095460", "ELVIS CANTILLO CASTILLA",            "2106", 129365, 3828128),
    @("73095460", "ELVIS CANTILLO CASTILLA",            "2107", 129365, 3828128),
    @("73095460", "ELVIS CANTILLO CASTILLA",            "2108", 129365, 3828128),
    @("73095460", "ELVIS CANTILLO CASTILLA",            "2109", 129365, 3828128),
    @("45481799", "LUZ ESTELA MARTINEZ LONDOÑO",        "2205", 240000, 7540320),
    @("45481799", "LUZ ESTELA MARTINEZ LONDOÑO",        "2206", 240000, 7540320),
    @("45481799", "LUZ ESTELA MARTINEZ LONDOÑO",        "2207", 240000, 7540320),
    @("45481799", "LUZ ESTELA MARTINEZ LONDOÑO",        "2208", 240000, 7540320),
    @("45481799", "LUZ ESTELA MARTINEZ LONDOÑO",        "2209", 240000, 7540320),
    @("45481799", "LUZ ESTELA MARTINEZ LONDOÑO",        "2210", 240000, 7540320),
    @("45481799", "LUZ ESTELA MARTINEZ LONDOÑO",        "2211", 240000, 7540320),
    @("23218869", "LEDYS DEL SOCORRO ROMERO PATERNINA",  "2409", 183674, 4591871),
    @("23218869", "LEDYS DEL SOCORRO ROMERO PATERNINA",  "2410", 183674, 4591871),
    @("23218869", "LEDYS DEL SOCORRO ROMERO PATERNINA",  "2411", 183674, 4591871),
    @("23218869", "LEDYS DEL SOCORRO ROMERO PATERNINA",  "2412", 183674, 4591871),
    @("23218869", "LEDYS DEL SOCORRO ROMERO PATERNINA",  "2501", 183674, 4591871),
    @("23218869", "LEDYS DEL SOCORRO ROMERO PATERNINA",  "2502", 183674, 4591871),
    @("23218869", "LEDYS DEL SOCORRO ROMERO PATERNINA",  "2503", 110205, 4591871)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 3).Value = $data[0]   # C - N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $data[1]   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $data[2]   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $data[3]   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $data[4]   # G - Salario Basico
}
